# Update "想去人数" (F column) counts on the "展览" and "全部类型" sheets
# to reflect newly scraped data (gh-pages output regenerated at 456a3b4).

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F5").Value = 11314
$ws1.Range("F9").Value = 11254
$ws1.Range("F11").Value = 1148
$ws1.Range("F14").Value = 5621
$ws1.Range("F15").Value = 102
$ws1.Range("F16").Value = 3465

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F7").Value = 11314
$ws4.Range("F11").Value = 11254
$ws4.Range("F13").Value = 1148
$ws4.Range("F16").Value = 5621
$ws4.Range("F17").Value = 102
$ws4.Range("F18").Value = 3465
